# Applies the crypto price/volume update described in the commit diff.
# Only cells whose content actually changes are touched; numeric-looking
# text values (e.g. "241.70", "1.000") are forced to Text format first
# so Excel does not silently convert them into numbers and drop information
# (trailing zeros, thousands-style dot grouping, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@(2, 4, '29.450.41', $false)
    ,@(2, 5, '  +0.09%  ', $false)
    ,@(3, 4, '1.857.34', $false)
    ,@(3, 5, '  +0.50%  ', $false)
    ,@(4, 4, '1.000', $true)
    ,@(4, 5, '  +0.01%  ', $false)
    ,@(5, 4, '241.70', $true)
    ,@(5, 5, '  +0.37%  ', $false)
    ,@(6, 4, '0.6339', $true)
    ,@(6, 5, '  +1.07%  ', $false)
    ,@(7, 4, '1.002', $true)
    ,@(7, 5, '  +0.11%  ', $false)
    ,@(8, 4, '0.07596', $true)
    ,@(8, 5, '  -1.10%  ', $false)
    ,@(9, 4, '0.2928', $true)
    ,@(9, 5, '  +0.40%  ', $false)
    ,@(10, 4, '24.61', $true)
    ,@(10, 5, '  -0.83%  ', $false)
    ,@(11, 4, '0.07767', $true)
    ,@(11, 5, '  +0.29%  ', $false)
    ,@(12, 4, '1.855.84', $false)
    ,@(12, 5, '  +0.57%  ', $false)
    ,@(13, 5, '  +0.15%  ', $false)
    ,@(14, 4, '0.6865', $true)
    ,@(14, 5, '  +0.93%  ', $false)
    ,@(15, 5, '  -2.69%  ', $false)
    ,@(16, 4, '83.43', $true)
    ,@(16, 5, '  +0.05%  ', $false)
    ,@(17, 2, 'WrappedliquidstakedEther2.0', $false)
    ,@(17, 3, 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', $false)
    ,@(17, 4, '2.113.83', $false)
    ,@(17, 5, '  +0.83%  ', $false)
    ,@(18, 2, 'Uniswap', $false)
    ,@(18, 3, 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', $false)
    ,@(18, 4, '6.160', $true)
    ,@(18, 5, '  -0.22%  ', $false)
    ,@(19, 2, 'WrappedBTC', $false)
    ,@(19, 3, 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', $false)
    ,@(19, 4, '29.445.68', $false)
    ,@(19, 5, '  +0.01%  ', $false)
    ,@(20, 2, 'BitcoinCash', $false)
    ,@(20, 3, 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', $false)
    ,@(20, 4, '230.69', $true)
    ,@(20, 5, '  +1.11%  ', $false)
    ,@(21, 2, 'Avalanche', $false)
    ,@(21, 3, 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', $false)
    ,@(21, 4, '12.40', $true)
    ,@(21, 5, '  +0.07%  ', $false)
    ,@(22, 2, 'Dai', $false)
    ,@(22, 3, 'https://coinranking.com/coin/MoTuySvg7+dai-dai', $false)
    ,@(22, 4, '1.001', $true)
    ,@(22, 5, '  +0.04%  ', $false)
    ,@(23, 2, 'Chainlink', $false)
    ,@(23, 3, 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', $false)
    ,@(23, 4, '7.535', $true)
    ,@(23, 5, '  +1.61%  ', $false)
    ,@(24, 2, 'BinanceUSD', $false)
    ,@(24, 3, 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', $false)
    ,@(24, 4, '1.001', $true)
    ,@(24, 5, '  -0.05%  ', $false)
    ,@(25, 2, 'Monero', $false)
    ,@(25, 3, 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', $false)
    ,@(25, 4, '159.19', $true)
    ,@(25, 5, '  +0.74%  ', $false)
    ,@(26, 2, 'Stellar', $false)
    ,@(26, 3, 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', $false)
    ,@(26, 4, '0.1404', $true)
    ,@(26, 5, '  +2.08%  ', $false)
    ,@(27, 2, 'Cosmos', $false)
    ,@(27, 3, 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', $false)
    ,@(27, 4, '8.485', $true)
    ,@(27, 5, '  +1.03%  ', $false)
    ,@(28, 2, 'EthereumClassic', $false)
    ,@(28, 3, 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', $false)
    ,@(28, 4, '17.77', $true)
    ,@(28, 5, '  +0.53%  ', $false)
    ,@(29, 2, 'Toncoin', $false)
    ,@(29, 3, 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', $false)
    ,@(29, 4, '1.419', $true)
    ,@(29, 5, '  +5.80%  ', $false)
    ,@(30, 2, 'PancakeSwap', $false)
    ,@(30, 3, 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', $false)
    ,@(30, 4, '1.481', $true)
    ,@(30, 5, '  +1.18%  ', $false)
    ,@(31, 2, 'Hedera', $false)
    ,@(31, 3, 'https://coinranking.com/coin/jad286TjB+hedera-hbar', $false)
    ,@(31, 4, '0.05706', $true)
    ,@(31, 5, '  +0.81%  ', $false)
    ,@(32, 2, 'Filecoin', $false)
    ,@(32, 3, 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', $false)
    ,@(32, 4, '4.161', $true)
    ,@(32, 5, '  +0.99%  ', $false)
    ,@(33, 2, 'InternetComputer(DFINITY)', $false)
    ,@(33, 3, 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', $false)
    ,@(33, 4, '4.066', $true)
    ,@(33, 5, '  +1.08%  ', $false)
    ,@(34, 2, 'LidoDAOToken', $false)
    ,@(34, 3, 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', $false)
    ,@(34, 4, '1.832', $true)
    ,@(34, 5, '  -0.48%  ', $false)
    ,@(35, 2, 'ARBITRUM', $false)
    ,@(35, 3, 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', $false)
    ,@(35, 4, '1.159', $true)
    ,@(35, 5, '  -0.23%  ', $false)
    ,@(36, 2, 'ImmutableX', $false)
    ,@(36, 3, 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', $false)
    ,@(36, 4, '0.6981', $true)
    ,@(36, 5, '  +0.55%  ', $false)
    ,@(37, 2, 'HuobiToken', $false)
    ,@(37, 3, 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', $false)
    ,@(37, 4, '2.590', $true)
    ,@(37, 5, '  +0.30%  ', $false)
    ,@(38, 2, 'Maker', $false)
    ,@(38, 3, 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', $false)
    ,@(38, 4, '1.251.26', $false)
    ,@(38, 5, '  +1.84%  ', $false)
    ,@(39, 2, 'VeChain', $false)
    ,@(39, 3, 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', $false)
    ,@(39, 4, '0.01828', $true)
    ,@(39, 5, '  +2.22%  ', $false)
    ,@(40, 2, 'MXToken', $false)
    ,@(40, 3, 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', $false)
    ,@(40, 4, '2.772', $true)
    ,@(40, 5, '  -0.01%  ', $false)
    ,@(41, 2, 'FraxShare', $false)
    ,@(41, 3, 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', $false)
    ,@(41, 4, '6.530', $true)
    ,@(41, 5, '  -0.09%  ', $false)
    ,@(42, 2, 'TrustWalletToken', $false)
    ,@(42, 3, 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', $false)
    ,@(42, 4, '0.9106', $true)
    ,@(42, 5, '  +0.16%  ', $false)
    ,@(43, 2, 'PaxDollar', $false)
    ,@(43, 3, 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', $false)
    ,@(43, 4, '1.001', $true)
    ,@(43, 5, '  +0.08%  ', $false)
    ,@(44, 2, 'RocketPoolETH', $false)
    ,@(44, 3, 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', $false)
    ,@(44, 4, '2.017.10', $false)
    ,@(44, 5, '  +0.60%  ', $false)
    ,@(45, 2, 'Quant', $false)
    ,@(45, 3, 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', $false)
    ,@(45, 4, '101.39', $true)
    ,@(45, 5, '  -0.34%  ', $false)
    ,@(46, 2, 'Aave', $false)
    ,@(46, 3, 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', $false)
    ,@(46, 4, '66.08', $true)
    ,@(46, 5, '  +0.21%  ', $false)
    ,@(47, 2, 'Aptos', $false)
    ,@(47, 3, 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', $false)
    ,@(47, 4, '7.172', $true)
    ,@(47, 5, '  +0.19%  ', $false)
    ,@(48, 2, 'Algorand', $false)
    ,@(48, 3, 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', $false)
    ,@(48, 4, '0.1167', $true)
    ,@(48, 5, '  +1.53%  ', $false)
    ,@(49, 2, 'EnergySwap', $false)
    ,@(49, 3, 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', $false)
    ,@(49, 4, '9.083', $true)
    ,@(49, 5, '  +0.91%  ', $false)
    ,@(50, 2, 'TheSandbox', $false)
    ,@(50, 3, 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', $false)
    ,@(50, 4, '0.3975', $true)
    ,@(50, 5, '  -1.00%  ', $false)
    ,@(51, 2, 'RenderToken', $false)
    ,@(51, 3, 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', $false)
    ,@(51, 4, '1.682', $true)
    ,@(51, 5, '  +0.56%  ', $false)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $forceText = $u[3]
    $cell = $ws.Cells.Item($row, $col)
    if ($forceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}
